# Add 7 new rows (44-50) of landscaping data below the existing table,
# extend the shared ABS(D-E) formula in column F, and update the
# active-cell selection to U2 (matching the post-edit state captured
# in the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New data rows (44-50). Column headers (row 1) are:
# A Date | B Plant_Type | C Plant_Size | D Low | E High | F Temp_Diff |
# G Rain | H Growth | I Pruned | J Quadrant | K Shade | L UV |
# M Humidity | N Dew_Point | O Pressure | P Wind_Gust | Q Cloud_Cover |
# R Visibility | S AQI | T Pollen
# ---------------------------------------------------------------------

$newRows = @(
    @{ Row=44; B="Flowering";     C="Large";  D=64; E=80; G=0.17; H=1;   I="No"; J=2; K="Bright";  L=7; M=0.67; N=67; O=29.74; P=26; Q=0.92; R=8.1; S=52; T=41 },
    @{ Row=45; B="Nonflowering";  C="Medium"; D=64; E=80; G=0.17; H=1;   I="No"; J=3; K="Bright";  L=7; M=0.67; N=67; O=29.74; P=26; Q=0.92; R=8.1; S=52; T=41 },
    @{ Row=46; B="Nonflowering";  C="Small";  D=64; E=80; G=0.17; H=1.5; I="No"; J=3; K="Neutral"; L=7; M=0.67; N=67; O=29.74; P=26; Q=0.92; R=8.1; S=52; T=41 },
    @{ Row=47; B="Nonflowering";  C="Medium"; D=64; E=80; G=0.17; H=0.5; I="No"; J=3; K="Dark";    L=7; M=0.67; N=67; O=29.74; P=26; Q=0.92; R=8.1; S=52; T=41 },
    @{ Row=48; B="Nonflowering";  C="Medium"; D=64; E=80; G=0.17; H=0.3; I="No"; J=3; K="Bright";  L=7; M=0.67; N=67; O=29.74; P=26; Q=0.92; R=8.1; S=52; T=41 },
    @{ Row=49; B="Nonflowering";  C="Large";  D=64; E=80; G=0.17; H=0.2; I="No"; J=4; K="Neutral"; L=7; M=0.67; N=67; O=29.74; P=26; Q=0.92; R=8.1; S=52; T=41 },
    @{ Row=50; B="Tree";          C="Medium"; D=64; E=80; G=0.17; H=0.1; I="No"; J=1; K="Neutral"; L=7; M=0.67; N=67; O=29.74; P=26; Q=0.92; R=8.1; S=52; T=41 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Column A: the date. Copy the format (style) from the last existing
    # data row (A43) so the new date cell reuses the existing date-number
    # style instead of creating a brand-new style entry.
    $ws.Cells.Item(43, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value = 45793

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    # Column F (Temp_Diff) is filled in afterwards as a shared formula.
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}

# Extend column F's ABS(Low-High difference) formula down through the
# newly added rows (F3:F43 -> F3:F50).
$ws.Range("F44:F50").Formula = "=ABS(D44-E44)"

# The worksheet's used range now runs through row 50 (dimension A1:U50
# is derived automatically from the populated cells).

# Update the current selection to match the saved state in the diff
# (activeCell U2, no frozen/scrolled topLeftCell).
$ws.Range("U2").Select()

# Best-effort: restore the workbook window size/position recorded in the
# diff. (Some hosts may not persist window geometry back to the saved
# file, since it is session/UI state rather than workbook content.)
$win = $excel.ActiveWindow
$win.Left = 672
$win.Top = 2844
$win.Width = 21156
$win.Height = 7128
